$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cell labels (B6, B8, B9) that identify which data cell feeds each
# confusion-matrix quantity (TP/FN/FP/TN).
$ws.Range("B6").Value = "C3"
$ws.Range("B8").Value = "C2"
$ws.Range("B9").Value = "B3"

# Fix the formulas that compute the derived statistics.
$ws.Range("B12").Formula = "=(100*C3)/(C3+C2)"
$ws.Range("B13").Formula = "=(100*B2)/(B2+B3)"
$ws.Range("B14").Formula = "=(100*(C3+B2))/(C3+B3+B2+C2)"
$ws.Range("B15").Formula = "=((C3*B2)-(B3*C2))/SQRT((C3+B3)*(C3+C2)*(B2+B3)*(B2+C2))"

# Move the active selection to D10, matching the saved view state.
$ws.Range("D10").Select()
